$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 573.5
$ws.Range("I28").Value = 528.06665
$ws.Range("J28").Value = 800.6667
$ws.Range("K28").Value = 528.06665
$ws.Range("L28").Value = 800.6667
$ws.Range("M28").Value = -43.06664999999998
$ws.Range("N28").Value = -1770.6667

$ws.Range("H93").Value = 24661.291
$ws.Range("J93").Value = 24661.291
$ws.Range("L93").Value = 24661.291
$ws.Range("N93").Value = -29653.291

$ws.Range("H113").Value = 4600.2383
$ws.Range("I113").Value = 2510.818
$ws.Range("J113").Value = 6898.6
$ws.Range("K113").Value = 2510.818
$ws.Range("L113").Value = 6898.6
$ws.Range("M113").Value = 743.1819999999998
$ws.Range("N113").Value = -13406.6

$ws.Range("H125").Value = 1575.3334
$ws.Range("J125").Value = 1575.3334
$ws.Range("L125").Value = 14178.0006
$ws.Range("N125").Value = -19098.0006

$ws.Range("H132").Value = 25004260
$ws.Range("I132").Value = 29415754
$ws.Range("J132").Value = 5794.3335
$ws.Range("K132").Value = 88247262
$ws.Range("L132").Value = 17383.0005
$ws.Range("M132").Value = -88244732
$ws.Range("N132").Value = -22443.0005

$ws.Range("H135").Value = 1268.5333
$ws.Range("I135").Value = 658.6667
$ws.Range("J135").Value = 2183.3333
$ws.Range("K135").Value = 5928.0003
$ws.Range("L135").Value = 19649.9997
$ws.Range("M135").Value = -3393.0003
$ws.Range("N135").Value = -24719.9997

$ws.Range("H137").Value = 1402479.1
$ws.Range("I137").Value = 1985703.8
$ws.Range("J137").Value = 2740.2
$ws.Range("K137").Value = 5957111.4
$ws.Range("L137").Value = 8220.599999999999
$ws.Range("M137").Value = -5954561.4
$ws.Range("N137").Value = -13320.6

$ws.Range("H138").Value = 2507.74
$ws.Range("I138").Value = 662.45
$ws.Range("J138").Value = 2969.0625
$ws.Range("K138").Value = 1987.35
$ws.Range("L138").Value = 8907.1875
$ws.Range("M138").Value = 3152.65
$ws.Range("N138").Value = -19187.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 702.6
$ws.Range("I2").Value = 650.375
$ws.Range("J2").Value = 911.5
$ws.Range("K2").Value = 650.375
$ws.Range("L2").Value = 911.5
$ws.Range("M2").Value = -537.375
$ws.Range("N2").Value = -1137.5

$ws.Range("H45").Value = 1390.4
$ws.Range("I45").Value = 1411.5555
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1411.5555
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -1034.5555
$ws.Range("N45").Value = -1954

$ws.Range("H61").Value = 1986.3572
$ws.Range("I61").Value = 1750
$ws.Range("J61").Value = 2004.5385
$ws.Range("K61").Value = 1750
$ws.Range("L61").Value = 2004.5385
$ws.Range("M61").Value = -1538
$ws.Range("N61").Value = -2428.5385

$ws.Range("H103").Value = 35063.43
$ws.Range("J103").Value = 35063.43
$ws.Range("L103").Value = 35063.43
$ws.Range("N103").Value = -37407.43

$ws.Range("H116").Value = 702.6
$ws.Range("I116").Value = 650.375
$ws.Range("J116").Value = 911.5
$ws.Range("K116").Value = 650.375
$ws.Range("L116").Value = 911.5
$ws.Range("M116").Value = 1643.625
$ws.Range("N116").Value = -5499.5

$ws.Range("H136").Value = 1986.3572
$ws.Range("I136").Value = 1750
$ws.Range("J136").Value = 2004.5385
$ws.Range("K136").Value = 5250
$ws.Range("L136").Value = 6013.6155
$ws.Range("M136").Value = -2700
$ws.Range("N136").Value = -11113.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 702.6
$ws.Range("I3").Value = 650.375
$ws.Range("J3").Value = 911.5
$ws.Range("K3").Value = 650.375
$ws.Range("L3").Value = 911.5
$ws.Range("M3").Value = -536.375
$ws.Range("N3").Value = -1139.5

$ws.Range("H95").Value = 30906.666
$ws.Range("J95").Value = 30906.666
$ws.Range("L95").Value = 30906.666
$ws.Range("N95").Value = -36398.666

$ws.Range("H103").Value = 38444.332
$ws.Range("J103").Value = 38444.332
$ws.Range("L103").Value = 38444.332
$ws.Range("N103").Value = -40788.332

$ws.Range("H134").Value = 2139.6956
$ws.Range("I134").Value = 1305.8235
$ws.Range("J134").Value = 4502.3335
$ws.Range("K134").Value = 3917.4705
$ws.Range("L134").Value = 13507.0005
$ws.Range("M134").Value = -1382.4705
$ws.Range("N134").Value = -18577.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2818.8965
$ws.Range("I31").Value = 1229.4546
$ws.Range("K31").Value = 1229.4546
$ws.Range("M31").Value = -934.4546

$ws.Range("H34").Value = 2818.8965
$ws.Range("I34").Value = 1229.4546
$ws.Range("K34").Value = 1229.4546
$ws.Range("M34").Value = -1027.4546

$ws.Range("H52").Value = 48600
$ws.Range("J52").Value = 48600
$ws.Range("L52").Value = 48600
$ws.Range("N52").Value = -49188

$ws.Range("H122").Value = 4254.1665
$ws.Range("I122").Value = 1855.5
$ws.Range("J122").Value = 5453.5
$ws.Range("K122").Value = 5566.5
$ws.Range("L122").Value = 16360.5
$ws.Range("M122").Value = -3116.5
$ws.Range("N122").Value = -21260.5

$ws.Range("H132").Value = 3579.303
$ws.Range("I132").Value = 3308.913
$ws.Range("J132").Value = 4201.2
$ws.Range("K132").Value = 9926.739
$ws.Range("L132").Value = 12603.6
$ws.Range("M132").Value = -7396.739
$ws.Range("N132").Value = -17663.6

$ws.Range("H134").Value = 6561.5835
$ws.Range("I134").Value = 8805.571
$ws.Range("J134").Value = 3420
$ws.Range("K134").Value = 26416.713
$ws.Range("L134").Value = 10260
$ws.Range("M134").Value = -23881.713
$ws.Range("N134").Value = -15330

$ws.Range("H137").Value = 41467.5
$ws.Range("J137").Value = 41467.5
$ws.Range("L137").Value = 41467.5
$ws.Range("N137").Value = -51667.5

$ws.Range("H139").Value = 38666.668
$ws.Range("J139").Value = 38666.668
$ws.Range("L139").Value = 38666.668
$ws.Range("N139").Value = -48946.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2925.4443

$ws.Range("H5").Value = 495551.88
$ws.Range("J5").Value = 891532.5600000001
$ws.Range("L5").Value = 2674597.68
$ws.Range("N5").Value = -2674821.68

$ws.Range("H113").Value = 4167262.2
$ws.Range("I113").Value = 602.2381
$ws.Range("K113").Value = 1806.7143
$ws.Range("M113").Value = 363.2856999999999

$ws.Range("H132").Value = 2310.5938
$ws.Range("J132").Value = 3020.9524
$ws.Range("L132").Value = 27188.5716
$ws.Range("N132").Value = -32248.5716

$ws.Range("H135").Value = 495551.88
$ws.Range("J135").Value = 891532.5600000001
$ws.Range("L135").Value = 8023793.040000001
$ws.Range("N135").Value = -8028863.040000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9618360
$ws.Range("I11").Value = 21400000
$ws.Range("J11").Value = 2254835.5
$ws.Range("K11").Value = 21400000
$ws.Range("L11").Value = 2254835.5
$ws.Range("M11").Value = -21399861
$ws.Range("N11").Value = -2255113.5

$ws.Range("H80").Value = 3387.4285
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 3118.6667
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 3118.6667
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -5114.6667

$ws.Range("H83").Value = 3387.4285
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 3118.6667
$ws.Range("K83").Value = 25000
$ws.Range("L83").Value = 15593.3335
$ws.Range("M83").Value = -20008
$ws.Range("N83").Value = -25577.3335

$ws.Range("H132").Value = 3554.7297
$ws.Range("I132").Value = 2840.4333
$ws.Range("K132").Value = 8521.2999
$ws.Range("M132").Value = -5991.2999

$ws.Range("H137").Value = 40506
$ws.Range("J137").Value = 40506
$ws.Range("L137").Value = 40506
$ws.Range("N137").Value = -50706

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4221.909
$ws.Range("I7").Value = 3746
$ws.Range("J7").Value = 5840
$ws.Range("K7").Value = 3746
$ws.Range("L7").Value = 5840
$ws.Range("M7").Value = -3634
$ws.Range("N7").Value = -6064

$ws.Range("H68").Value = 725.85
$ws.Range("I68").Value = 723.1818
$ws.Range("K68").Value = 723.1818
$ws.Range("M68").Value = 25.81820000000005

$ws.Range("H71").Value = 725.85
$ws.Range("I71").Value = 723.1818
$ws.Range("K71").Value = 3615.909
$ws.Range("M71").Value = 128.0910000000003

$ws.Range("H126").Value = 4221.909
$ws.Range("I126").Value = 3746
$ws.Range("J126").Value = 5840
$ws.Range("K126").Value = 11238
$ws.Range("L126").Value = 17520
$ws.Range("M126").Value = -8768
$ws.Range("N126").Value = -22460

$ws.Range("H132").Value = 6504.727
$ws.Range("I132").Value = 2847.077
$ws.Range("J132").Value = 11788
$ws.Range("K132").Value = 8541.231
$ws.Range("L132").Value = 35364
$ws.Range("M132").Value = -6011.231
$ws.Range("N132").Value = -40424

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13589.286
$ws.Range("I136").Value = 13783.375
$ws.Range("J136").Value = 13330.5
$ws.Range("K136").Value = 41350.125
$ws.Range("L136").Value = 39991.5
$ws.Range("M136").Value = -38800.125
$ws.Range("N136").Value = -45091.5
